$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values need to be swapped between row 13 and row 14
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range($col + "13")
    $cell14 = $ws.Range($col + "14")

    $v13 = $cell13.Value2
    $v14 = $cell14.Value2

    $cell13.Value2 = $v14
    $cell14.Value2 = $v13
}
